$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 25 (item #24, profile.php / "Ielogojoties ar lietotaju..."): mark as resolved ---
# Copy cell formatting from row 23 (an existing "resolved" row with the same column style pattern)
# so the fill/border/alignment match the "A - atrisinats" group exactly, then update status + resolver.
$ws.Range("A23:G23").Copy() | Out-Null
$ws.Range("A25:G25").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(25, 5).Value = "A"
$ws.Cells.Item(25, 6).Value = "Arnis"

# --- Row 27 (item #26, searchRoom.php / "Auditorijas adrese"): trim resolved note, shrink row ---
$ws.Cells.Item(27, 4).Value = '1. Jāpieliek lauks "Auditorijas adrese"'
$ws.Rows.Item(27).AutoFit() | Out-Null

# --- Row 28 (item #27, searchCourse.php / "Divreiz atkartojas..."): mark as resolved ---
$ws.Range("A23:G23").Copy() | Out-Null
$ws.Range("A28:G28").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(28, 5).Value = "A"
$ws.Cells.Item(28, 6).Value = "Arnis"

# --- Row 29 (item #28, newRoom.php / "Studentu skaits..."): mark as resolved ---
$ws.Range("A23:G23").Copy() | Out-Null
$ws.Range("A29:G29").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(29, 5).Value = "A"
$ws.Cells.Item(29, 6).Value = "Arnis"

$excel.CutCopyMode = 0

# --- Move the saved view position / selection to the new working area ---
$ws.Range("F35").Select() | Out-Null
